$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.590.51"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").Value = "2.004.50"
$ws.Range("E3").Value = "  -3.98%  "
$ws.Range("E4").Value = "  +0.91%  "
$ws.Range("D5").Value = "'329.42"
$ws.Range("E5").Value = "  -3.87%  "
$ws.Range("D6").Value = "'1.010"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").Value = "'0.5015"
$ws.Range("E7").Value = "  -4.04%  "
$ws.Range("D8").Value = "'0.4229"
$ws.Range("E8").Value = "  -3.88%  "
$ws.Range("D9").Value = "'54.08"
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("D10").Value = "'0.09033"
$ws.Range("E10").Value = "  -3.25%  "
$ws.Range("D11").Value = "'1.120"
$ws.Range("E11").Value = "  -3.92%  "
$ws.Range("D12").Value = "'23.35"
$ws.Range("E12").Value = "  -5.57%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'8.075"
$ws.Range("E13").Value = "  -5.79%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.979.16"
$ws.Range("E14").Value = "  -3.24%  "
$ws.Range("D15").Value = "'6.496"
$ws.Range("E15").Value = "  -5.62%  "
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "'94.28"
$ws.Range("E17").Value = "  -7.04%  "
$ws.Range("D18").Value = "'0.00001114"
$ws.Range("E18").Value = "  -3.60%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "'19.70"
$ws.Range("E20").Value = "  -6.54%  "
$ws.Range("D21").Value = "'1.011"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").Value = "'5.968"
$ws.Range("E22").Value = "  -5.54%  "
$ws.Range("D23").Value = "29.634.80"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("D24").Value = "'12.03"
$ws.Range("E24").Value = "  -3.59%  "
$ws.Range("D25").Value = "'2.300"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").Value = "'158.80"
$ws.Range("E26").Value = "  -2.07%  "
$ws.Range("E27").Value = "  -4.70%  "
$ws.Range("D28").Value = "'6.400"
$ws.Range("E28").Value = "  -3.79%  "
$ws.Range("E29").Value = "  -7.67%  "
$ws.Range("D30").Value = "'128.21"
$ws.Range("E30").Value = "  -3.53%  "
$ws.Range("E31").Value = "  -6.15%  "
$ws.Range("D32").Value = "'0.09972"
$ws.Range("E32").Value = "  -4.55%  "
$ws.Range("D33").Value = "'1.572"
$ws.Range("E33").Value = "  -4.98%  "
$ws.Range("D34").Value = "'5.849"
$ws.Range("E34").Value = "  -5.82%  "
$ws.Range("D35").Value = "'3.794"
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("D36").Value = "'0.02477"
$ws.Range("E36").Value = "  -5.61%  "
$ws.Range("D37").Value = "'9.325"
$ws.Range("E37").Value = "  -8.42%  "
$ws.Range("D38").Value = "'1.312"
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("D39").Value = "'0.06381"
$ws.Range("E39").Value = "  -6.09%  "
$ws.Range("D40").Value = "'0.6570"
$ws.Range("E40").Value = "  -5.47%  "
$ws.Range("E41").Value = "  -6.11%  "
$ws.Range("D42").Value = "'0.2056"
$ws.Range("E42").Value = "  -6.70%  "
$ws.Range("D43").Value = "'1.011"
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").Value = "'0.6357"
$ws.Range("E44").Value = "  -6.29%  "
$ws.Range("D45").Value = "'13.43"
$ws.Range("E45").Value = "  -5.80%  "
$ws.Range("D46").Value = "'2.197"
$ws.Range("E46").Value = "  -5.40%  "
$ws.Range("D47").Value = "'1.306"
$ws.Range("E47").Value = "  -4.66%  "
$ws.Range("D48").Value = "'3.510"
$ws.Range("E48").Value = "  -3.34%  "
$ws.Range("D49").Value = "'0.00000000336"
$ws.Range("E49").Value = "  -4.44%  "
$ws.Range("D50").Value = "'0.06989"
$ws.Range("E50").Value = "  -3.25%  "
$ws.Range("E51").Value = "  -6.70%  "
